$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '47.784.66'
$ws.Range('D2').NumberFormat = 'General'
$ws.Range('E2').Value = '  +5.52%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.624.46'
$ws.Range('D3').NumberFormat = 'General'
$ws.Range('E3').Value = '  +9.81%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').NumberFormat = 'General'
$ws.Range('E4').Value = '  +0.05%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '313.02'
$ws.Range('D5').NumberFormat = 'General'
$ws.Range('E5').Value = '  +6.62%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '105.43'
$ws.Range('D6').NumberFormat = 'General'
$ws.Range('E6').Value = '  +12.42%  '

$ws.Range('E7').Value = '  +9.41%  '

$ws.Range('E9').Value = '  +18.61%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '40.03'
$ws.Range('D10').NumberFormat = 'General'
$ws.Range('E10').Value = '  +16.97%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '55.31'
$ws.Range('D11').NumberFormat = 'General'
$ws.Range('E11').Value = '  +3.65%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0858'
$ws.Range('D12').NumberFormat = 'General'
$ws.Range('E12').Value = '  +10.39%  '

$ws.Range('E13').Value = '  +19.40%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '3.028.36'
$ws.Range('D14').NumberFormat = 'General'
$ws.Range('E14').Value = '  +9.88%  '

$ws.Range('E15').Value = '  +3.19%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.628.27'
$ws.Range('D16').NumberFormat = 'General'
$ws.Range('E16').Value = '  +10.07%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.938'
$ws.Range('D17').NumberFormat = 'General'
$ws.Range('E17').Value = '  +13.79%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '15.24'
$ws.Range('D18').NumberFormat = 'General'
$ws.Range('E18').Value = '  +8.60%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '47.791.32'
$ws.Range('D19').NumberFormat = 'General'
$ws.Range('E19').Value = '  +5.63%  '

$ws.Range('E20').Value = '  +10.45%  '

$ws.Range('E21').Value = '  +5.79%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.77'
$ws.Range('D22').NumberFormat = 'General'
$ws.Range('E22').Value = '  +11.17%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '72.95'
$ws.Range('D23').NumberFormat = 'General'
$ws.Range('E23').Value = '  +9.63%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '271.98'
$ws.Range('D24').NumberFormat = 'General'
$ws.Range('E24').Value = '  +14.14%  '

$ws.Range('E25').Value = '  +12.21%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.23'
$ws.Range('D26').NumberFormat = 'General'
$ws.Range('E26').Value = '  +18.07%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '30.50'
$ws.Range('D27').NumberFormat = 'General'
$ws.Range('E27').Value = '  +45.68%  '

$ws.Range('E28').Value = '  +0.15%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '4.10'
$ws.Range('D29').NumberFormat = 'General'
$ws.Range('E29').Value = '  +1.65%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '10.68'
$ws.Range('D30').NumberFormat = 'General'
$ws.Range('E30').Value = '  +11.86%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '39.90'
$ws.Range('D31').NumberFormat = 'General'
$ws.Range('E31').Value = '  +6.76%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.30'
$ws.Range('D32').NumberFormat = 'General'
$ws.Range('E32').Value = '  +4.03%  '

$ws.Range('E33').Value = '  +14.20%  '

$ws.Range('E34').Value = '  -2.96%  '

$ws.Range('E35').Value = '  +7.13%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0854'
$ws.Range('D36').NumberFormat = 'General'
$ws.Range('E36').Value = '  +12.61%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.21'
$ws.Range('D37').NumberFormat = 'General'
$ws.Range('E37').Value = '  +12.98%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '152.19'
$ws.Range('D38').NumberFormat = 'General'
$ws.Range('E38').Value = '  +3.20%  '

$ws.Range('E39').Value = '  +11.47%  '

$ws.Range('E40').Value = '  +8.72%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '16.38'
$ws.Range('D41').NumberFormat = 'General'
$ws.Range('E41').Value = '  +12.93%  '

$ws.Range('B42').Value = 'RenderToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '4.29'
$ws.Range('D42').NumberFormat = 'General'
$ws.Range('E42').Value = '  +16.05%  '

$ws.Range('B43').Value = 'EnergySwap'
$ws.Range('C43').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '22.85'
$ws.Range('D43').NumberFormat = 'General'
$ws.Range('E43').Value = '  +52.72%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.68'
$ws.Range('D44').NumberFormat = 'General'
$ws.Range('E44').Value = '  +16.77%  '

$ws.Range('E45').Value = '  +14.01%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.200.96'
$ws.Range('D46').NumberFormat = 'General'
$ws.Range('E46').Value = '  +11.87%  '

$ws.Range('B47').Value = 'FirstDigitalUSD'
$ws.Range('C47').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.999'
$ws.Range('D47').NumberFormat = 'General'
$ws.Range('E47').Value = '  +0.08%  '

$ws.Range('B48').Value = 'BitcoinSV'
$ws.Range('C48').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '95.44'
$ws.Range('D48').NumberFormat = 'General'
$ws.Range('E48').Value = '  +8.01%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '10.06'
$ws.Range('D49').NumberFormat = 'General'
$ws.Range('E49').Value = '  +19.35%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '113.98'
$ws.Range('D50').NumberFormat = 'General'
$ws.Range('E50').Value = '  +14.93%  '

$ws.Range('E51').Value = '  +5.91%  '
